# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (column F) and one ticket-price
# (column G) change across the 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibition) sheet ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 481
$ws1.Range("F4").Value = 368
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F5").Value = 1787
$ws1.Range("F6").Value = 391
$ws1.Range("F7").Value = 1465
$ws1.Range("F9").Value = 366
$ws1.Range("F10").Value = 713
$ws1.Range("F11").Value = 13038
$ws1.Range("F12").Value = 12955
$ws1.Range("F13").Value = 974
$ws1.Range("F14").Value = 755
$ws1.Range("F18").Value = 614
$ws1.Range("F19").Value = 2035
$ws1.Range("F21").Value = 21
$ws1.Range("F22").Value = 24
$ws1.Range("F24").Value = 152
$ws1.Range("F26").Value = 718

# ---- 演出 (Performance) sheet ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 25
$ws2.Range("F7").Value = 92
$ws2.Range("F10").Value = 7

# ---- 全部类型 (All types) sheet ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 481
$ws4.Range("F6").Value = 368
$ws4.Range("G6").Value = "不可售"
$ws4.Range("F7").Value = 1787
$ws4.Range("F8").Value = 391
$ws4.Range("F9").Value = 1465
$ws4.Range("F11").Value = 366
$ws4.Range("F13").Value = 713
$ws4.Range("F14").Value = 13038
$ws4.Range("F15").Value = 12955
$ws4.Range("F16").Value = 974
$ws4.Range("F17").Value = 755
$ws4.Range("F21").Value = 614
$ws4.Range("F22").Value = 25
$ws4.Range("F24").Value = 2035
$ws4.Range("F26").Value = 21
$ws4.Range("F27").Value = 24
$ws4.Range("F31").Value = 152
$ws4.Range("F33").Value = 718
$ws4.Range("F34").Value = 92
$ws4.Range("F37").Value = 7
